$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.104.59'
$ws.Range("E2").Value = '  -0.19%  '
$ws.Range("D3").Value = '3.476.27'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.67'
$ws.Range("E5").Value = '  -0.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.89'
$ws.Range("E6").Value = '  -1.63%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -1.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.71'
$ws.Range("E9").Value = '  +6.03%  '
$ws.Range("E10").Value = '  -1.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.386'
$ws.Range("E11").Value = '  +0.13%  '
$ws.Range("D12").Value = '4.072.31'
$ws.Range("E12").Value = '  -0.28%  '
$ws.Range("E13").Value = '  +0.05%  '
$ws.Range("E14").Value = '  -2.51%  '
$ws.Range("D15").Value = '3.479.84'
$ws.Range("E15").Value = '  -0.28%  '
$ws.Range("D16").Value = '64.121.49'
$ws.Range("E16").Value = '  -0.29%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.11'
$ws.Range("E17").Value = '  -2.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '9.97'
$ws.Range("E18").Value = '  +1.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.69'
$ws.Range("E19").Value = '  -1.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.38'
$ws.Range("E20").Value = '  -1.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '384.75'
$ws.Range("E21").Value = '  -2.36%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.566'
$ws.Range("E22").Value = '  -0.43%  '
$ws.Range("D23").Value = '3.619.42'
$ws.Range("E23").Value = '  -0.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.48'
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("B26").Value = 'PEPE'
$ws.Range("C26").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000111'
$ws.Range("E26").Value = '  -3.51%  '
$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  -0.23%  '
$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.22'
$ws.Range("E28").Value = '  -0.57%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.08'
$ws.Range("E29").Value = '  -4.18%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.98'
$ws.Range("E30").Value = '  -3.11%  '
$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.155'
$ws.Range("E31").Value = '  +2.98%  '
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.42'
$ws.Range("E32").Value = '  -4.61%  '
$ws.Range("B33").Value = 'RenzoRestakedETH'
$ws.Range("C33").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D33").Value = '3.507.13'
$ws.Range("E33").Value = '  -0.14%  '
$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  -0.10%  '
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '22.96'
$ws.Range("E35").Value = '  -1.83%  '
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.21'
$ws.Range("E36").Value = '  +1.42%  '
$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.76'
$ws.Range("E37").Value = '  -2.02%  '
$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '162.36'
$ws.Range("E38").Value = '  -2.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.49'
$ws.Range("E39").Value = '  -3.56%  '
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0777'
$ws.Range("E40").Value = '  -0.34%  '
$ws.Range("B41").Value = 'Mantle'
$ws.Range("C41").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.798'
$ws.Range("E41").Value = '  -0.89%  '
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.51'
$ws.Range("E43").Value = '  -0.78%  '
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.33'
$ws.Range("E44").Value = '  -1.20%  '
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.62'
$ws.Range("E45").Value = '  -1.90%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.52'
$ws.Range("E46").Value = '  -7.09%  '
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.13'
$ws.Range("E47").Value = '  -3.43%  '
$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.71'
$ws.Range("E48").Value = '  -0.85%  '
$ws.Range("B49").Value = 'SuiNetwork'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.899'
$ws.Range("E49").Value = '  +0.39%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '2.336.72'
$ws.Range("E50").Value = '  -4.92%  '
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0253'
$ws.Range("E51").Value = '  -2.69%  '
